$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/22/2025  Through  12/28/2025"

# --- Row 16-21, 23-31 numeric updates ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 141
$ws.Range("J16").Value = 140
$ws.Range("K16").Value = 0.714285714285
$ws.Range("L16").Value = -18.023255813953
$ws.Range("M16").Value = -4.729729729729
$ws.Range("N16").Value = -84.175084175084
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 47.058823529411
$ws.Range("I17").Value = 296
$ws.Range("J17").Value = 268
$ws.Range("K17").Value = 10.447761194029
$ws.Range("L17").Value = 16.99604743083
$ws.Range("M17").Value = 166.666666666667
$ws.Range("N17").Value = -5.431309904153
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -76.190476190476
$ws.Range("I18").Value = 143
$ws.Range("J18").Value = 176
$ws.Range("K18").Value = -18.75
$ws.Range("L18").Value = -35.874439461883
$ws.Range("M18").Value = -42.570281124498
$ws.Range("N18").Value = -93.273753527751
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 60
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = -7.5
$ws.Range("I19").Value = 569
$ws.Range("J19").Value = 588
$ws.Range("K19").Value = -3.231292517006
$ws.Range("L19").Value = -11.919504643962
$ws.Range("M19").Value = 42.606516290726
$ws.Range("N19").Value = -60.839642119752
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -60.714285714285
$ws.Range("I20").Value = 225
$ws.Range("J20").Value = 273
$ws.Range("K20").Value = -17.582417582417
$ws.Range("L20").Value = 11.940298507462
$ws.Range("M20").Value = 21.621621621621
$ws.Range("N20").Value = -93.125572868927
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 15
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 122
$ws.Range("H21").Value = -29.508196721311
$ws.Range("I21").Value = 1391
$ws.Range("J21").Value = 1464
$ws.Range("K21").Value = -4.986338797814
$ws.Range("L21").Value = -8.124174372523
$ws.Range("M21").Value = 26.22504537205
$ws.Range("N21").Value = -82.799554841103
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 62
$ws.Range("J23").Value = 77
$ws.Range("K23").Value = -19.480519480519
$ws.Range("L23").Value = 1.639344262295
$ws.Range("M23").Value = 106.666666666667
$ws.Range("C24").Value = 19
$ws.Range("E24").Value = -17.391304347826
$ws.Range("F24").Value = 79
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = -21.782178217821
$ws.Range("I24").Value = 1114
$ws.Range("J24").Value = 1202
$ws.Range("K24").Value = -7.321131447587
$ws.Range("L24").Value = -11.022364217252
$ws.Range("M24").Value = 2.578268876611
$ws.Range("C25").Value = 7
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = -28.205128205128
$ws.Range("I25").Value = 409
$ws.Range("J25").Value = 505
$ws.Range("K25").Value = -19.009900990099
$ws.Range("L25").Value = -22.975517890772
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 71.428571428571
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = -9.090909090909
$ws.Range("I26").Value = 519
$ws.Range("J26").Value = 497
$ws.Range("K26").Value = 4.426559356136
$ws.Range("L26").Value = 25.970873786407
$ws.Range("M26").Value = 3.386454183266
$ws.Range("D28").Value = 1
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = 44.444444444444
$ws.Range("D29").Value = 1
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -66.666666666666
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = -44.444444444444
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -86.486486486486
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = -50
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -16.666666666666
$ws.Range("M30").Value = -28.571428571428
$ws.Range("N30").Value = -83.333333333333
$ws.Range("G31").Value = 1
$ws.Range("L31").Value = -5.555555555555

# --- Row 22: D22/E22 switch from numeric to text placeholder cells ---
# Copy the formatting+value from existing text-placeholder cells in the same
# row so the style index and shared-string text match exactly.
$ws.Range("C22").Copy($ws.Range("D22"))
$ws.Range("N22").Copy($ws.Range("E22"))

# --- Row 23: C23 switches from text placeholder "0" to numeric 1 ---
$ws.Range("D23").Copy($ws.Range("C23"))
